$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 28-31 need to be filled in with the same sample-row pattern used by
# rows 2-27 (bioSampleNumber/rnaSampleNumber continue incrementing).
# Copying the previous, already-correctly-typed/styled row across is the
# most reliable way to replicate the exact cell types (shared-string text
# vs. number) and styles (e.g. the existing "s=1" text style already
# present on column H) without Excel's automatic value-sniffing turning
# the date-like strings into real dates or "False" into a boolean.

$ws.Range("A27:H27").Copy($ws.Range("A28:H28"))
$ws.Range("C28").Value = 27
$ws.Range("F28").Value = 27

$ws.Range("A28:H28").Copy($ws.Range("A29:H29"))
$ws.Range("C29").Value = 28
$ws.Range("F29").Value = 28

$ws.Range("A29:H29").Copy($ws.Range("A30:H30"))
$ws.Range("C30").Value = 29
$ws.Range("F30").Value = 29

$ws.Range("A30:H30").Copy($ws.Range("A31:H31"))
$ws.Range("C31").Value = 30
$ws.Range("F31").Value = 30

# Match the author's final cursor position recorded in the saved sheet view.
$ws.Range("I22").Select()
